$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# New pseudo-random sample data for the color-scale conditional formatting demo
# (the conditional formatting rules on B2:H5, D7:G12 and C15:J18 stay untouched;
# only the underlying A1:J20 values are refreshed).
$values = @(
    @(55.8603827170378, 46.3780275761979, 152.243937995398, 78.9956533717903, 135.261113538994, 46.9787431168271, 91.0547565161505, 108.087243655737, 131.641277918425, 26.7630994444541),
    @(141.496570101705, 178.752195638955, 1.22404526976126, 199.819334689444, 149.854011624052, 105.145647239473, 199.634231254288, 35.8542567285962, 39.9836600944324, 172.899967885064),
    @(171.821368379435, 183.284240580762, 99.8080708551258, 179.59213768113, 96.4630233573089, 5.09410631148802, 187.341186212069, 132.477039998619, 107.713346233458, 40.3243849241754),
    @(70.8726690480824, 28.3576615286794, 69.8911891644314, 117.031572161723, 107.564645590058, 66.6489888292965, 137.187778547959, 144.921344399881, 29.2703065226182, 152.582909330997),
    @(78.0191712444737, 183.04097791344, 188.052196888277, 98.3287659000274, 55.9140092022317, 11.3105726480999, 60.8798746303096, 13.4657994906724, 194.316295345461, 181.430495707984),
    @(3.9219872112954, 163.538239041128, 25.0210194964991, 33.8525410899206, 181.53562768434, 72.5761421362758, 146.569956721072, 172.651800314268, 182.532630014481, 130.167007227506),
    @(59.6375569047581, 158.577716517531, 0.373897422279183, 91.3168929942497, 155.890430396372, 113.138908573025, 108.861006474523, 84.1924731080385, 92.2546890993857, 83.2050227947556),
    @(167.957868691514, 54.7128868544069, 6.58395020597798, 87.4007507634353, 94.8807966405902, 188.780390465995, 195.232043692485, 1.47930495509845, 123.678128478899, 85.1524507092091),
    @(144.214231681178, 173.875386721397, 138.160744653158, 126.282850525474, 36.40239771288, 107.334430006954, 3.33664203218028, 36.0386480745108, 135.495944477383, 34.9885034537821),
    @(120.079032108224, 164.535978233691, 162.3887143854, 99.1032992951122, 92.945352426239, 119.441454726943, 182.667080491161, 96.7353038940278, 142.438335503656, 142.775100629206),
    @(102.449566173577, 176.687401522271, 121.211110391287, 111.111272550705, 13.4726270164701, 149.209100356889, 156.95428883515, 137.620268733064, 138.97174444933, 192.755237218344),
    @(77.3440984437913, 145.090651765974, 48.9736718353693, 97.3801793052723, 185.952775546328, 85.7621701833616, 20.4169718643729, 74.0910468968056, 54.9144952813697, 48.5560003894176),
    @(109.802266540845, 72.8223584000125, 148.696528630656, 57.2661856456037, 163.125990686531, 3.42189045782289, 92.3241724690069, 107.480650910866, 194.455398337196, 175.439341913648),
    @(6.11330997483493, 98.4967397984568, 59.0409694514428, 180.903027849692, 182.702884535633, 167.526830158907, 52.6642763301098, 27.0494721024528, 112.810223509004, 87.1932973559915),
    @(150.380141171804, 65.7163732991164, 97.0669036251804, 142.740707259039, 157.644405009991, 174.988380342251, 115.562306398322, 65.0085350801277, 113.150523748785, 7.1831822428774),
    @(99.0244828625696, 108.576033594355, 41.8208086126581, 93.882335114238, 32.9728340883613, 29.627207773564, 27.9908728916156, 63.9449247456831, 147.985281864174, 10.0507365586472),
    @(56.8849278878816, 49.4736379242845, 143.164870395868, 163.532402535683, 186.641927243509, 178.847358645335, 86.0496823145308, 68.0706439856769, 114.67729476964, 18.4259453874202),
    @(33.0978938532518, 193.36749976192, 161.280823387802, 167.721197925378, 98.1758592176139, 44.0858932417286, 175.755454774832, 5.95582137161671, 99.6217806356129, 188.13761034428),
    @(87.8595840595009, 27.3156373888793, 194.330127162081, 187.272216094319, 76.414859051078, 97.5372763804799, 56.6759311857987, 165.158634337205, 147.930193761331, 153.075676762068),
    @(139.535957267292, 188.719351584427, 79.0641902382784, 102.759486950356, 30.3083694681099, 100.906503247519, 122.551502903249, 133.534501089498, 156.09878001553, 178.797046364656),
)

for ($r = 1; $r -le 20; $r++) {
    for ($c = 1; $c -le 10; $c++) {
        $ws.Cells.Item($r, $c).Value = $values[$r - 1][$c - 1]
    }
}

Write-Output "Updated A1:J20 with new sample values"
